$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 220 (i.e. at positions 221-222),
# pushing the former rows 221..316 down to 223..318.
$ws.Rows.Item(221).Resize(2).Insert()

# --- Fill in the first new row (221) ---
$ws.Range("A221").Value = 5
$ws.Range("B221").Value = "Macroferia Regional de Talca"
$ws.Range("C221").Value = "Maule"
$ws.Range("D221").Value = 44460
$ws.Range("E221").Value = 7
$ws.Range("F221").Value = 100112002
$ws.Range("G221").Value = "Pimiento"
$ws.Range("H221").Value = "Zafiro rojo"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 200
$ws.Range("K221").Value = 35000
$ws.Range("L221").Value = 35000
$ws.Range("M221").Value = 35000
$ws.Range("N221").Value = "`$/caja 15 kilos"
$ws.Range("O221").Value = "Región de Arica y Parinacota"
$ws.Range("P221").Value = 2333
$ws.Range("Q221").Value = 15
$ws.Range("R221").Value = "Hortaliza"

# --- Fill in the second new row (222) ---
$ws.Range("A222").Value = 5
$ws.Range("B222").Value = "Macroferia Regional de Talca"
$ws.Range("C222").Value = "Maule"
$ws.Range("D222").Value = 44460
$ws.Range("E222").Value = 7
$ws.Range("F222").Value = 100112002
$ws.Range("G222").Value = "Pimiento"
$ws.Range("H222").Value = "Zafiro verde"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 200
$ws.Range("K222").Value = 25000
$ws.Range("L222").Value = 25000
$ws.Range("M222").Value = 25000
$ws.Range("N222").Value = "`$/caja 15 kilos"
$ws.Range("O222").Value = "Región de Arica y Parinacota"
$ws.Range("P222").Value = 1667
$ws.Range("Q222").Value = 15
$ws.Range("R222").Value = "Hortaliza"
